# Edit: insert 4 new weekly price rows for "Hortaliza, Vega Monumental Concepción - Cebolla"
# right after the existing row 251 (i.e. they become the new rows 252-255),
# pushing the previously existing rows 252-276 down to 256-280.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before current row 252 (they inherit formatting from row 251 above).
$ws.Rows("252:255").Insert()

# Common / repeated values for this subset (same market/region/category across the sheet)
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$categoriaId = 100112004
$categoria = "Cebolla"
$clasificacion = "Hortaliza"

# New row 252
$r = 252
$ws.Cells.Item($r, 1).Value() = $mercadoId
$ws.Cells.Item($r, 2).Value() = $mercado
$ws.Cells.Item($r, 3).Value() = $region
$ws.Cells.Item($r, 4).Value() = [DateTime]"2021-10-22"
$ws.Cells.Item($r, 5).Value() = $codreg
$ws.Cells.Item($r, 6).Value() = $categoriaId
$ws.Cells.Item($r, 7).Value() = $categoria
$ws.Cells.Item($r, 8).Value() = "Morada(o)"
$ws.Cells.Item($r, 9).Value() = "1a (cosecha)"
$ws.Cells.Item($r, 10).Value() = 200
$ws.Cells.Item($r, 11).Value() = 6000
$ws.Cells.Item($r, 12).Value() = 6500
$ws.Cells.Item($r, 13).Value() = 6250
$ws.Cells.Item($r, 14).Value() = "$/malla 18 kilos"
$ws.Cells.Item($r, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item($r, 16).Value() = 347
$ws.Cells.Item($r, 17).Value() = 18
$ws.Cells.Item($r, 18).Value() = $clasificacion

# New row 253
$r = 253
$ws.Cells.Item($r, 1).Value() = $mercadoId
$ws.Cells.Item($r, 2).Value() = $mercado
$ws.Cells.Item($r, 3).Value() = $region
$ws.Cells.Item($r, 4).Value() = [DateTime]"2021-10-22"
$ws.Cells.Item($r, 5).Value() = $codreg
$ws.Cells.Item($r, 6).Value() = $categoriaId
$ws.Cells.Item($r, 7).Value() = $categoria
$ws.Cells.Item($r, 8).Value() = "Morada(o)"
$ws.Cells.Item($r, 9).Value() = "2a (cosecha)"
$ws.Cells.Item($r, 10).Value() = 100
$ws.Cells.Item($r, 11).Value() = 5500
$ws.Cells.Item($r, 12).Value() = 5500
$ws.Cells.Item($r, 13).Value() = 5500
$ws.Cells.Item($r, 14).Value() = "$/malla 18 kilos"
$ws.Cells.Item($r, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item($r, 16).Value() = 306
$ws.Cells.Item($r, 17).Value() = 18
$ws.Cells.Item($r, 18).Value() = $clasificacion

# New row 254
$r = 254
$ws.Cells.Item($r, 1).Value() = $mercadoId
$ws.Cells.Item($r, 2).Value() = $mercado
$ws.Cells.Item($r, 3).Value() = $region
$ws.Cells.Item($r, 4).Value() = [DateTime]"2021-10-22"
$ws.Cells.Item($r, 5).Value() = $codreg
$ws.Cells.Item($r, 6).Value() = $categoriaId
$ws.Cells.Item($r, 7).Value() = $categoria
$ws.Cells.Item($r, 8).Value() = "Sin especificar"
$ws.Cells.Item($r, 9).Value() = "1a nueva(o)"
$ws.Cells.Item($r, 10).Value() = 1000
$ws.Cells.Item($r, 11).Value() = 1800
$ws.Cells.Item($r, 12).Value() = 2000
$ws.Cells.Item($r, 13).Value() = 1900
$ws.Cells.Item($r, 14).Value() = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item($r, 15).Value() = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value() = 95
$ws.Cells.Item($r, 17).Value() = 20
$ws.Cells.Item($r, 18).Value() = $clasificacion

# New row 255
$r = 255
$ws.Cells.Item($r, 1).Value() = $mercadoId
$ws.Cells.Item($r, 2).Value() = $mercado
$ws.Cells.Item($r, 3).Value() = $region
$ws.Cells.Item($r, 4).Value() = [DateTime]"2021-10-22"
$ws.Cells.Item($r, 5).Value() = $codreg
$ws.Cells.Item($r, 6).Value() = $categoriaId
$ws.Cells.Item($r, 7).Value() = $categoria
$ws.Cells.Item($r, 8).Value() = "Sin especificar"
$ws.Cells.Item($r, 9).Value() = "2a nueva(o)"
$ws.Cells.Item($r, 10).Value() = 500
$ws.Cells.Item($r, 11).Value() = 1600
$ws.Cells.Item($r, 12).Value() = 1600
$ws.Cells.Item($r, 13).Value() = 1600
$ws.Cells.Item($r, 14).Value() = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item($r, 15).Value() = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value() = 80
$ws.Cells.Item($r, 17).Value() = 20
$ws.Cells.Item($r, 18).Value() = $clasificacion

# Make sure date column keeps the same date-only number format as the rest of column D
$ws.Range("D252:D255").NumberFormat = $ws.Range("D256").NumberFormat()
